$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")

# Mark the six "cyan spheres / track end baubles" tasks (rows 86-91, column C) as done
$ws.Range("C86:C91").Value = "X"

# Move the active selection to where the cursor would land after filling these cells
$ws.Activate()
$ws.Range("C92").Select()
